$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CageNumber 1): height 20 -> 25, Material Metal -> Plastic
$ws.Range("C2").Value = 25
$ws.Range("E2").Value = "Plastic"

# Row 4 (CageNumber 3): length 60 -> 30
$ws.Range("B4").Value = 30

# Row 5 (CageNumber 4): height 30 -> 15, width 40 -> 20, Material Wood -> Plastic
$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = "Plastic"

# Row 9: CageNumber 8 -> 9, height 15 -> 30, width 20 -> 40, Material Plastic -> Wood
$ws.Range("A9").Value = 9
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 40
$ws.Range("E9").Value = "Wood"

# Row 10: CageNumber 9 -> 10, length 20 -> 15, width 40 -> 25, Material Wood -> Plastic
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = 15
$ws.Range("D10").Value = 25
$ws.Range("E10").Value = "Plastic"

# Row 11: CageNumber 10 -> 11, length 15 -> 25, width 25 -> 32, Material Plastic -> Wood
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = 25
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = "Wood"

# Row 12 (new cage): CageNumber 12, length 25, height 30, width 40, Material Wood
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = 25
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 40
$ws.Range("E12").Value = "Wood"

# Row 13 (new cage): CageNumber 15, length 30, height 20, width 20, Material Plastic
$ws.Range("A13").Value = 15
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 20
$ws.Range("E13").Value = "Plastic"

Write-Host "Edits applied"
